$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, D contain text that can look numeric (e.g. "0.9980",
# "27.595.40"); force Text format before assigning so Excel keeps the
# exact original string (no numeric coercion / trailing-zero loss),
# then restore the default "Normal" style so no stray formatting is left
# behind (Column E values already contain non-numeric characters like
# '%' / surrounding spaces, so they stay text without the extra step).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.595.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.775.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9953'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4438'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.66%  '
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.06'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07727'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.119'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9980'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.73'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.182'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.442'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.765.60'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.41'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +11.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001074'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06281'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9966'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.207'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.5293'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.640.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.40%  '
$ws.Range("E25").Value = '  -2.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.264'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.317'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.967.65'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '128.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.184'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.66%  '
$ws.Range("E33").Value = '  -2.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09198'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.656'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.67'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02317'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2167'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06127'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6472'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.072'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.183'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.950'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9958'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.386'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5987'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.717'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '126.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.990'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.83%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Cronos'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06911'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.65%  '
